$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new day's price row was scraped (2025-11-24). Insert it above the prior
# top data row, pushing the existing rows (2025-11-23 .. 2025-11-21) down by
# one, and carry the latest metric readings into the new row.
$ws.Rows("2:2").Insert()

# Insert() copies the formatting of the row above (the bold/bordered header),
# so strip that back off the freshly inserted row before writing data into it.
$ws.Range("A2:D2").ClearFormats()

# Force column A to stay plain text so "2025-11-24" isn't reinterpreted as a
# date serial number, matching how the date column is stored elsewhere in
# the sheet.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-11-24"
$ws.Range("A2").ClearFormats()

$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
